$d = $word.ActiveDocument

# --- Text replacements in the API table ---
$d.Content.Find.Execute("INTER", $true, $false, $false, $false, $false, $true, 1, $false, "ROUT", 2)
$d.Content.Find.Execute("PROV_INT", $true, $false, $false, $false, $false, $true, 1, $false, "PROV_ROUT", 2)
$d.Content.Find.Execute("SET_INT", $true, $false, $false, $false, $false, $true, 1, $false, "SET_ROUT", 2)
$d.Content.Find.Execute("DEL_INT", $true, $false, $false, $false, $false, $true, 1, $false, "DEL_ROUT", 2)

Write-Host "Done"
